# Update the expiration-date test data on the eCard test-cases sheet
# and move the active selection, per the source commit
# ":memo: TEST DATA - added new test data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eCard_testovacie karty")

# D4 / D5 hold the "Datum expiracie" (expiration date) values for the two
# test cards; bump them from 2024-12-31 to 2028-12-31 (serial 47118).
$ws.Range("D4").Value = 47118
$ws.Range("D5").Value = 47118

# Reflect the new active cell/selection left by the editor.
$ws.Activate()
$ws.Range("D12").Select()
